# Updated CHE_grids model - 2025-08-09 17:26
# Re-assign the "grid_cell" label (column AG, rows 4-26) on the "solar"
# sheet's elc_won/elc_spv bus-connection lookup table to a new ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

$gridCells = @{
    4  = "CHE_24"
    5  = "CHE_5"
    6  = "CHE_8"
    7  = "CHE_13"
    8  = "CHE_21"
    9  = "CHE_9"
    10 = "CHE_4"
    11 = "CHE_14"
    12 = "CHE_18"
    13 = "CHE_3"
    14 = "CHE_20"
    15 = "CHE_1"
    16 = "CHE_6"
    17 = "CHE_17"
    18 = "CHE_19"
    19 = "CHE_10"
    20 = "CHE_22"
    21 = "CHE_11"
    22 = "CHE_15"
    23 = "CHE_25"
    24 = "CHE_0"
    25 = "CHE_7"
    26 = "CHE_12"
}

foreach ($row in $gridCells.Keys) {
    $ws.Range("AG$row").Value = $gridCells[$row]
}
